$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E7").Value = 18
$ws.Range("E9").Value = 14
$ws.Range("E10").Value = 14
$ws.Range("E13").Value = 3
$ws.Range("E14").Value = 28
$ws.Range("E15").Value = 66
$ws.Range("E16").Value = 232
$ws.Range("E18").Value = 64
